$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 105, shifting rows 105:179 down to 106:180
$ws.Rows.Item(105).Insert()

# Populate the new row 105 with the new data point
$ws.Cells.Item(105, 1).Value = 7
$ws.Cells.Item(105, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(105, 3).Value = "Ñuble"
$ws.Cells.Item(105, 4).Value = 44574
$ws.Cells.Item(105, 5).Value = 16
$ws.Cells.Item(105, 6).Value = 100112003
$ws.Cells.Item(105, 7).Value = "Ajo"
$ws.Cells.Item(105, 8).Value = "Chino"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 100
$ws.Cells.Item(105, 11).Value = 20000
$ws.Cells.Item(105, 12).Value = 21000
$ws.Cells.Item(105, 13).Value = 20500
$ws.Cells.Item(105, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(105, 15).Value = "China"
$ws.Cells.Item(105, 16).Value = 2050
$ws.Cells.Item(105, 17).Value = 10
$ws.Cells.Item(105, 18).Value = "Hortaliza"
